# Saldo_guide.xlsx update: refresh the "Dt. Referencia" reference date by one
# day (2024-04-23 -> 2024-04-24, serial 45405 -> 45406) for every data row,
# and correct three balances whose "Saldo Previsto"/"Vl. Total" entries
# (columns D and H) were overstated by an extra leading 1/16/12 thousand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to the last used row (row 310 in this workbook).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 7)
    if ($dateCell.Value2 -ne $null) {
        $dateCell.Value = $dateCell.Value2 + 1
    }
}

# Balance corrections (Saldo Previsto = column D, Vl. Total = column H)
# Row 2:   84601.76 -> 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 8).Value = 0

# Row 55:  16205.42 -> 205.42
$ws.Cells.Item(55, 4).Value = 205.42
$ws.Cells.Item(55, 8).Value = 205.42

# Row 135: 12234.67 -> 234.67
$ws.Cells.Item(135, 4).Value = 234.67
$ws.Cells.Item(135, 8).Value = 234.67
